# Weekly update: add a new week's data (2 new quality-grade rows) at the top
# of the data block (row 442), pushing all existing data rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 442.
$ws.Rows("442:443").Insert()

# New row 442 - "Primera" quality for the new reporting date.
$ws.Cells(442, 1).Value2 = 4
$ws.Cells(442, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells(442, 3).Value2 = "Los Lagos"
$ws.Cells(442, 4).Value2 = 45127
$ws.Cells(442, 5).Value2 = 10
$ws.Cells(442, 6).Value2 = 100112017
$ws.Cells(442, 7).Value2 = "Apio"
$ws.Cells(442, 8).Value2 = "Americana (o)"
$ws.Cells(442, 9).Value2 = "Primera"
$ws.Cells(442, 10).Value2 = 15
$ws.Cells(442, 11).Value2 = 11000
$ws.Cells(442, 12).Value2 = 11000
$ws.Cells(442, 13).Value2 = 11000
$ws.Cells(442, 14).Value2 = "$/docena de matas"
$ws.Cells(442, 15).Value2 = "Región de Coquimbo"
$ws.Cells(442, 16).Value2 = 1833
$ws.Cells(442, 17).Value2 = 6
$ws.Cells(442, 18).Value2 = "Hortaliza"

# New row 443 - "Segunda" quality for the same new reporting date.
$ws.Cells(443, 1).Value2 = 4
$ws.Cells(443, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells(443, 3).Value2 = "Los Lagos"
$ws.Cells(443, 4).Value2 = 45127
$ws.Cells(443, 5).Value2 = 10
$ws.Cells(443, 6).Value2 = 100112017
$ws.Cells(443, 7).Value2 = "Apio"
$ws.Cells(443, 8).Value2 = "Americana (o)"
$ws.Cells(443, 9).Value2 = "Segunda"
$ws.Cells(443, 10).Value2 = 15
$ws.Cells(443, 11).Value2 = 10000
$ws.Cells(443, 12).Value2 = 10000
$ws.Cells(443, 13).Value2 = 10000
$ws.Cells(443, 14).Value2 = "$/docena de matas"
$ws.Cells(443, 15).Value2 = "Región de Coquimbo"
$ws.Cells(443, 16).Value2 = 1667
$ws.Cells(443, 17).Value2 = 6
$ws.Cells(443, 18).Value2 = "Hortaliza"
